# Add an "AVERAGE" column (E) to the measurement sheets: 6kV, d=10cm cambio
# volt, and V1=6V d=10cm. Column E header = "AVERAGE"; cells below hold
# =AVERAGE(Bn:Dn).

$wb = $excel.ActiveWorkbook

# --- Sheet "6kV": rows 2-8, average column already pre-formatted (style s=2)
$ws = $wb.Worksheets.Item("6kV")
$ws.Range("E1").Value = "AVERAGE"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=AVERAGE(B${r}:D${r})"
}

# --- Sheet "d=10cm cambio volt": rows 2-6, average column needs "0.0" format
$ws = $wb.Worksheets.Item("d=10cm cambio volt")
$ws.Range("E1").Value = "AVERAGE"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=AVERAGE(B${r}:D${r})"
}
$ws.Range("E2:E6").NumberFormat = "0.0"

# --- Sheet "V1=6V d=10cm": rows 2-6, average column needs "0.0" format
$ws = $wb.Worksheets.Item("V1=6V d=10cm")
$ws.Range("E1").Value = "AVERAGE"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=AVERAGE(B${r}:D${r})"
}
$ws.Range("E2:E6").NumberFormat = "0.0"
